$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 550
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = -174
$ws.Range("N32").Value = -1252
$ws.Range("H42").Value = 174.5
$ws.Range("J42").Value = 236.85715
$ws.Range("L42").Value = 710.5714499999999
$ws.Range("N42").Value = -1170.57145
$ws.Range("H98").Value = 466299.97
$ws.Range("I98").Value = 620955.75
$ws.Range("J98").Value = 2332.6667
$ws.Range("K98").Value = 620955.75
$ws.Range("L98").Value = 2332.6667
$ws.Range("M98").Value = -619457.75
$ws.Range("N98").Value = -5328.6667
$ws.Range("H122").Value = 466299.97
$ws.Range("I122").Value = 620955.75
$ws.Range("J122").Value = 2332.6667
$ws.Range("K122").Value = 1862867.25
$ws.Range("L122").Value = 6998.000100000001
$ws.Range("M122").Value = -1860417.25
$ws.Range("N122").Value = -11898.0001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21104.607
$ws.Range("I32").Value = 2766.423
$ws.Range("K32").Value = 2766.423
$ws.Range("M32").Value = -2479.423

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2257.8
$ws.Range("I20").Value = 2125.182
$ws.Range("J20").Value = 2419.889
$ws.Range("K20").Value = 2125.182
$ws.Range("L20").Value = 2419.889
$ws.Range("M20").Value = -1878.182
$ws.Range("N20").Value = -2913.889
$ws.Range("H86").Value = 6872.6
$ws.Range("I86").Value = 2214.5
$ws.Range("J86").Value = 11530.7
$ws.Range("K86").Value = 2214.5
$ws.Range("L86").Value = 11530.7
$ws.Range("M86").Value = -1091.5
$ws.Range("N86").Value = -13776.7
$ws.Range("H89").Value = 6872.6
$ws.Range("I89").Value = 2214.5
$ws.Range("J89").Value = 11530.7
$ws.Range("K89").Value = 11072.5
$ws.Range("L89").Value = 57653.5
$ws.Range("M89").Value = -5456.5
$ws.Range("N89").Value = -68885.5
$ws.Range("H134").Value = 34485740
$ws.Range("I134").Value = 47620772
$ws.Range("K134").Value = 142862316
$ws.Range("M134").Value = -142859781

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3832.8667
$ws.Range("I31").Value = 1028
$ws.Range("J31").Value = 16805.375
$ws.Range("K31").Value = 1028
$ws.Range("L31").Value = 16805.375
$ws.Range("M31").Value = -733
$ws.Range("N31").Value = -17395.375
$ws.Range("H34").Value = 3832.8667
$ws.Range("I34").Value = 1028
$ws.Range("J34").Value = 16805.375
$ws.Range("K34").Value = 1028
$ws.Range("L34").Value = 16805.375
$ws.Range("M34").Value = -826
$ws.Range("N34").Value = -17209.375
$ws.Range("H122").Value = 1509.7727
$ws.Range("I122").Value = 1276.25
$ws.Range("J122").Value = 1790
$ws.Range("K122").Value = 3828.75
$ws.Range("L122").Value = 5370
$ws.Range("M122").Value = -1378.75
$ws.Range("N122").Value = -10270
$ws.Range("H132").Value = 3050.6667
$ws.Range("I132").Value = 2627.8518
$ws.Range("J132").Value = 4953.3335
$ws.Range("K132").Value = 7883.555399999999
$ws.Range("L132").Value = 14860.0005
$ws.Range("M132").Value = -5353.555399999999
$ws.Range("N132").Value = -19920.0005
$ws.Range("H134").Value = 2947.4827
$ws.Range("I134").Value = 1419.95
$ws.Range("J134").Value = 6342
$ws.Range("K134").Value = 4259.85
$ws.Range("L134").Value = 19026
$ws.Range("M134").Value = -1724.85
$ws.Range("N134").Value = -24096
$ws.Range("H135").Value = 38715
$ws.Range("J135").Value = 38715
$ws.Range("L135").Value = 38715
$ws.Range("N135").Value = -48855

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1604.75
$ws.Range("I18").Value = 139.66667
$ws.Range("J18").Value = 6000
$ws.Range("K18").Value = 419.00001
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = -250.00001
$ws.Range("N18").Value = -18338
$ws.Range("H51").Value = 1309.7142
$ws.Range("I51").Value = 3252
$ws.Range("K51").Value = 9756
$ws.Range("M51").Value = -9296
$ws.Range("H55").Value = 2477.25
$ws.Range("I55").Value = 1004
$ws.Range("J55").Value = 2968.3333
$ws.Range("K55").Value = 3012
$ws.Range("L55").Value = 8904.999899999999
$ws.Range("M55").Value = -2835
$ws.Range("N55").Value = -9258.999899999999
$ws.Range("H59").Value = 1000
$ws.Range("I59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -2460
$ws.Range("H61").Value = 362.55554
$ws.Range("I61").Value = 305.6
$ws.Range("J61").Value = 433.75
$ws.Range("K61").Value = 916.8000000000001
$ws.Range("L61").Value = 1301.25
$ws.Range("M61").Value = -701.8000000000001
$ws.Range("N61").Value = -1731.25
$ws.Range("H80").Value = 1027.5714
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1048.8334
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3146.5002
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -5018.5002
$ws.Range("H83").Value = 1027.5714
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1048.8334
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 9439.500599999999
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -18799.5006
$ws.Range("H113").Value = 664.2162
$ws.Range("I113").Value = 648
$ws.Range("J113").Value = 683.2941
$ws.Range("K113").Value = 1944
$ws.Range("L113").Value = 2049.8823
$ws.Range("M113").Value = 226
$ws.Range("N113").Value = -6389.882299999999
$ws.Range("H141").Value = 3755.0833
$ws.Range("I141").Value = 4056.1
$ws.Range("J141").Value = 2250
$ws.Range("K141").Value = 12168.3
$ws.Range("L141").Value = 6750
$ws.Range("M141").Value = -6988.299999999999
$ws.Range("N141").Value = -17110

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2607
$ws.Range("I80").Value = 2470.5557
$ws.Range("J80").Value = 2852.6
$ws.Range("K80").Value = 2470.5557
$ws.Range("L80").Value = 2852.6
$ws.Range("M80").Value = -1472.5557
$ws.Range("N80").Value = -4848.6
$ws.Range("H83").Value = 2607
$ws.Range("I83").Value = 2470.5557
$ws.Range("J83").Value = 2852.6
$ws.Range("K83").Value = 12352.7785
$ws.Range("L83").Value = 14263
$ws.Range("M83").Value = -7360.7785
$ws.Range("N83").Value = -24247
$ws.Range("H122").Value = 531267.3
$ws.Range("I122").Value = 856368.3
$ws.Range("J122").Value = 2978.25
$ws.Range("K122").Value = 2569104.9
$ws.Range("L122").Value = 8934.75
$ws.Range("M122").Value = -2566654.9
$ws.Range("N122").Value = -13834.75
$ws.Range("H126").Value = 2566
$ws.Range("I126").Value = 1996.6666
$ws.Range("J126").Value = 2708.3333
$ws.Range("K126").Value = 5989.9998
$ws.Range("L126").Value = 8124.999899999999
$ws.Range("M126").Value = -3519.9998
$ws.Range("N126").Value = -13064.9999
$ws.Range("H132").Value = 4084.35
$ws.Range("I132").Value = 3812.4
$ws.Range("J132").Value = 4356.3
$ws.Range("K132").Value = 11437.2
$ws.Range("L132").Value = 13068.9
$ws.Range("M132").Value = -8907.200000000001
$ws.Range("N132").Value = -18128.9
$ws.Range("H133").Value = 44550
$ws.Range("J133").Value = 44550
$ws.Range("L133").Value = 44550
$ws.Range("N133").Value = -54670

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 30007
$ws.Range("J43").Value = 10014
$ws.Range("L43").Value = 10014
$ws.Range("N43").Value = -10400
$ws.Range("H55").Value = 818.2
$ws.Range("I55").Value = 899.5
$ws.Range("J55").Value = 797.875
$ws.Range("K55").Value = 899.5
$ws.Range("L55").Value = 797.875
$ws.Range("M55").Value = -726.5
$ws.Range("N55").Value = -1143.875
$ws.Range("H82").Value = 2698.9092
$ws.Range("I82").Value = 2460.6667
$ws.Range("J82").Value = 2788.25
$ws.Range("K82").Value = 2460.6667
$ws.Range("L82").Value = 2788.25
$ws.Range("M82").Value = -2099.6667
$ws.Range("N82").Value = -3510.25
$ws.Range("H85").Value = 2698.9092
$ws.Range("I85").Value = 2460.6667
$ws.Range("J85").Value = 2788.25
$ws.Range("K85").Value = 2460.6667
$ws.Range("L85").Value = 2788.25
$ws.Range("M85").Value = -1212.6667
$ws.Range("N85").Value = -5284.25
$ws.Range("H122").Value = 2690.25
$ws.Range("I122").Value = 1820.421
$ws.Range("J122").Value = 3961.5386
$ws.Range("K122").Value = 5461.263
$ws.Range("L122").Value = 11884.6158
$ws.Range("M122").Value = -3011.263
$ws.Range("N122").Value = -16784.6158

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 167834
$ws.Range("I122").Value = 200800.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 602402.3999999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -599952.3999999999
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 100578.2
$ws.Range("I126").Value = 100578.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 301734.6
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -299264.6
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2447.4075
$ws.Range("I136").Value = 1250
$ws.Range("K136").Value = 3750
$ws.Range("M136").Value = -1200
